$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "44.653.42"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +1.56%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.243.68"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +0.44%  "

# Row 4
$ws.Range("E4").Value = "  +0.24%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.10"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.70%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "96.17"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +2.22%  "

# Row 7
$ws.Range("E7").Value = "  +0.32%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.07%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.523"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.49%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.20"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +1.52%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0806"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.00%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.26"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +1.33%  "

# Row 13
$ws.Range("E13").Value = "  +0.18%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.300.41"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +3.01%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.839"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +2.46%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.64"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +1.19%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "44.349.14"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +1.27%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0₃0954"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.49%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.34"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.45%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.02"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.42%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "65.60"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.42%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "238.31"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.88%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.97"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +2.44%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.99"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +2.33%  "

# Row 25
$ws.Range("E25").Value = "  +0.10%  "

# Row 26
$ws.Range("E26").Value = "  +4.47%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.84"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.02%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "37.22"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -1.31%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.00"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.34%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "20.01"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.74%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "152.58"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +1.12%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0800"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.61%  "

# Row 33
$ws.Range("E33").Value = "  +2.46%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.06"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -3.87%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.109"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.11%  "

# Row 36
$ws.Range("E36").Value = "  -0.17%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.84"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +4.43%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "15.05"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +1.27%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.41"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +2.35%  "

# Row 40
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0301"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +2.45%  "

# Row 41
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.75"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -1.70%  "

# Row 42
$ws.Range("E42").Value = "  +0.17%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.824.00"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +5.66%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.71"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +15.64%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.194"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +4.32%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "79.40"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -5.98%  "

# Row 47
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "99.33"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.02%  "

# Row 48
$ws.Range("B48").Value = "ordi"
$ws.Range("C48").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "70.95"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +3.79%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.92"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.28%  "

# Row 50
$ws.Range("B50").Value = "MultiversX"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "54.75"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +2.13%  "

# Row 51
$ws.Range("B51").Value = "FraxShare"
$ws.Range("C51").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.06"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.41%  "
